$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.473.58"
$ws.Range("E2").Value = "  -1.92%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.832.71"
$ws.Range("E3").Value = "  -2.69%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.60%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.69"
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4595"
$ws.Range("E7").Value = "  -3.79%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3818"
$ws.Range("E8").Value = "  -3.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.60"
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07915"
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9710"
$ws.Range("E11").Value = "  -4.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.08"
$ws.Range("E12").Value = "  -4.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.829.49"
$ws.Range("E13").Value = "  -3.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.891"
$ws.Range("E14").Value = "  -2.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.027"
$ws.Range("E15").Value = "  -2.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.006"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.00"
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06652"
$ws.Range("E18").Value = "  -1.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001029"
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.00"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.469.53"
$ws.Range("E22").Value = "  -1.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.348"
$ws.Range("E23").Value = "  -3.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.82"
$ws.Range("E24").Value = "  -1.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.308"
$ws.Range("E25").Value = "  -1.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.051.91"
$ws.Range("E26").Value = "  -2.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.41"
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.37"
$ws.Range("E28").Value = "  -2.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.059"
$ws.Range("E29").Value = "  -2.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.297"
$ws.Range("E30").Value = "  -3.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.01"
$ws.Range("E31").Value = "  -2.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9539"
$ws.Range("E32").Value = "  -2.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09296"
$ws.Range("E33").Value = "  -3.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.590"
$ws.Range("E34").Value = "  -1.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.247"
$ws.Range("E35").Value = "  -1.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.312"
$ws.Range("E36").Value = "  -3.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05933"
$ws.Range("E37").Value = "  -2.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02192"
$ws.Range("E38").Value = "  -2.65%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.163"
$ws.Range("E39").Value = "  -3.93%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.063"
$ws.Range("E40").Value = "  -1.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5787"
$ws.Range("E41").Value = "  -3.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1838"
$ws.Range("E42").Value = "  -3.38%  "
$ws.Range("E43").Value = "  -3.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.253"
$ws.Range("E44").Value = "  -1.03%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5481"
$ws.Range("E45").Value = "  -3.76%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.97"
$ws.Range("E46").Value = "  -1.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.863"
$ws.Range("E47").Value = "  -3.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06646"
$ws.Range("E48").Value = "  -2.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.43"
$ws.Range("E49").Value = "  -2.22%  "
$ws.Range("E50").Value = "  -3.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.005"
$ws.Range("E51").Value = "  -0.66%  "
